# chore: simulator full-month coverage, persist logs, fix employees
#
# - Timesheet entries for the week standardized to 8-hr "Regular" days
#   at a $140/day rate (was a mix of partial hours / PTO at $0 rate).
# - Several client names on the timesheet corrected (simulator data fix).
# - Employee ID corrected for Chris Zavesky.
# - The "Jason Schema" log sheet mirrors the same corrections and now
#   persists Type / Notes consistently on every row (row 2 previously
#   leaked "PTO" into the Type/Notes columns instead of Regular/blank).

$wb  = $excel.ActiveWorkbook
$tsw = $wb.Worksheets.Item("Weekly Timesheet")
$log = $wb.Worksheets.Item("Jason Schema")

# ---------------------------------------------------------------------
# Weekly Timesheet (columns: A Date | B Client | C Hours | D Type |
#                             E Rate | F Total)
# ---------------------------------------------------------------------

# Row 2 - 2026-01-19
$tsw.Range("B2").Value = "Tormey"
$tsw.Range("C2").Value = 8
$tsw.Range("D2").Value = "Regular"
$tsw.Range("E2").Value = 140
$tsw.Range("F2").Value = 1120

# Row 3 - 2026-01-20
$tsw.Range("B3").Value = "Evans"
$tsw.Range("C3").Value = 8
$tsw.Range("D3").Value = "Regular"
$tsw.Range("E3").Value = 140
$tsw.Range("F3").Value = 1120

# Row 4 - 2026-01-21
$tsw.Range("C4").Value = 8
$tsw.Range("D4").Value = "Regular"
$tsw.Range("E4").Value = 140
$tsw.Range("F4").Value = 1120

# Row 5 - 2026-01-22
$tsw.Range("B5").Value = "Hewett"
$tsw.Range("C5").Value = 8
$tsw.Range("D5").Value = "Regular"
$tsw.Range("E5").Value = 140
$tsw.Range("F5").Value = 1120

# Row 6 - 2026-01-23
$tsw.Range("B6").Value = "Howard"
$tsw.Range("C6").Value = 8
$tsw.Range("D6").Value = "Regular"
$tsw.Range("E6").Value = 140
$tsw.Range("F6").Value = 1120

# Row 8 - SUBTOTAL (Reg hours + label + $ total)
$tsw.Range("C8").Value = 40
$tsw.Range("D8").Value = "Reg: 40 / OT: 0"
$tsw.Range("F8").Value = 5600

# Row 12 - HOURLY SUBTOTAL ($ total)
$tsw.Range("F12").Value = 5600

# Row 13 - GRAND TOTAL ($ total)
$tsw.Range("F13").Value = 5600

# ---------------------------------------------------------------------
# Jason Schema (columns: A Employee | B Employee ID | C Date | D Client |
#                         E Hours | F Rate | G Total | H Type | I Notes)
# ---------------------------------------------------------------------

# Row 2 - 2026-01-19
$log.Range("D2").Value = "Tormey"
$log.Range("E2").Value = 8
$log.Range("F2").Value = 140
$log.Range("G2").Value = 1120
$log.Range("H2").Value = "Regular"
$log.Range("I2").Value = ""

# Row 3 - 2026-01-20
$log.Range("D3").Value = "Evans"
$log.Range("E3").Value = 8
$log.Range("F3").Value = 140
$log.Range("G3").Value = 1120
$log.Range("H3").Value = "Regular"

# Row 4 - 2026-01-21
$log.Range("E4").Value = 8
$log.Range("F4").Value = 140
$log.Range("G4").Value = 1120
$log.Range("H4").Value = "Regular"

# Row 5 - 2026-01-22
$log.Range("D5").Value = "Hewett"
$log.Range("E5").Value = 8
$log.Range("F5").Value = 140
$log.Range("G5").Value = 1120
$log.Range("H5").Value = "Regular"

# Row 6 - 2026-01-23
$log.Range("D6").Value = "Howard"
$log.Range("E6").Value = 8
$log.Range("F6").Value = 140
$log.Range("G6").Value = 1120
$log.Range("H6").Value = "Regular"

# ---------------------------------------------------------------------
# Employee ID correction - appears on both sheets
# ---------------------------------------------------------------------
$oldId = "emp_5chpvt65"
$newId = "emp_jp4mlvog"

for ($r = 2; $r -le 6; $r++) {
    if ($log.Cells.Item($r, 2).Value2 -eq $oldId) {
        $log.Cells.Item($r, 2).Value = $newId
    }
}
